$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Characters().Text = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.206.06"
Set-TextValue $ws.Range("E2") "  +2.25%  "
Set-TextValue $ws.Range("D3") "1.906.44"
Set-TextValue $ws.Range("E3") "  +2.09%  "
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("D5") "327.86"
Set-TextValue $ws.Range("E5") "  +1.08%  "
Set-TextValue $ws.Range("D6") "1.004"
Set-TextValue $ws.Range("E6") "  +0.10%  "
Set-TextValue $ws.Range("D7") "0.4660"
Set-TextValue $ws.Range("E7") "  +1.15%  "
Set-TextValue $ws.Range("D8") "0.3972"
Set-TextValue $ws.Range("E8") "  +2.67%  "
Set-TextValue $ws.Range("D9") "46.82"
Set-TextValue $ws.Range("E9") "  +1.12%  "
Set-TextValue $ws.Range("D10") "0.07969"
Set-TextValue $ws.Range("E10") "  +1.47%  "
Set-TextValue $ws.Range("D11") "1.003"
Set-TextValue $ws.Range("E11") "  +3.06%  "
Set-TextValue $ws.Range("D12") "22.35"
Set-TextValue $ws.Range("E12") "  +2.13%  "
Set-TextValue $ws.Range("D13") "1.909.45"
Set-TextValue $ws.Range("E13") "  +2.90%  "
Set-TextValue $ws.Range("D14") "7.149"
Set-TextValue $ws.Range("E14") "  +2.54%  "
Set-TextValue $ws.Range("D15") "5.792"
Set-TextValue $ws.Range("E15") "  +1.78%  "
Set-TextValue $ws.Range("D16") "0.06950"
Set-TextValue $ws.Range("E16") "  +0.22%  "
Set-TextValue $ws.Range("D17") "88.88"
Set-TextValue $ws.Range("E17") "  +0.87%  "
Set-TextValue $ws.Range("E18") "  +0.12%  "
Set-TextValue $ws.Range("D19") "0.00001013"
Set-TextValue $ws.Range("E19") "  +1.14%  "
Set-TextValue $ws.Range("D20") "17.20"
Set-TextValue $ws.Range("E20") "  +2.48%  "
Set-TextValue $ws.Range("D21") "1.002"
Set-TextValue $ws.Range("E21") "  -0.05%  "
Set-TextValue $ws.Range("D22") "29.223.09"
Set-TextValue $ws.Range("E22") "  +2.31%  "
Set-TextValue $ws.Range("D23") "5.358"
Set-TextValue $ws.Range("E23") "  +1.83%  "
Set-TextValue $ws.Range("D24") "11.09"
Set-TextValue $ws.Range("E24") "  +0.63%  "
Set-TextValue $ws.Range("D25") "2.130.72"
Set-TextValue $ws.Range("E25") "  +2.75%  "
Set-TextValue $ws.Range("D26") "2.059"
Set-TextValue $ws.Range("E26") "  -2.41%  "
Set-TextValue $ws.Range("D27") "156.75"
Set-TextValue $ws.Range("E27") "  +2.87%  "
Set-TextValue $ws.Range("D28") "19.53"
Set-TextValue $ws.Range("E28") "  +1.62%  "
Set-TextValue $ws.Range("D29") "5.887"
Set-TextValue $ws.Range("E29") "  +2.08%  "
Set-TextValue $ws.Range("E30") "  +0.80%  "
Set-TextValue $ws.Range("D31") "119.61"
Set-TextValue $ws.Range("E31") "  +0.37%  "
Set-TextValue $ws.Range("D32") "0.09443"
Set-TextValue $ws.Range("E32") "  +1.24%  "
Set-TextValue $ws.Range("D33") "0.9234"
Set-TextValue $ws.Range("E33") "  +0.65%  "
Set-TextValue $ws.Range("D34") "5.356"
Set-TextValue $ws.Range("E34") "  +1.88%  "
Set-TextValue $ws.Range("D35") "1.346"
Set-TextValue $ws.Range("E35") "  +1.08%  "
Set-TextValue $ws.Range("D36") "3.265"
Set-TextValue $ws.Range("E36") "  -1.84%  "
Set-TextValue $ws.Range("D37") "0.05852"
Set-TextValue $ws.Range("E37") "  +1.19%  "
Set-TextValue $ws.Range("D38") "1.175"
Set-TextValue $ws.Range("E38") "  +1.84%  "
Set-TextValue $ws.Range("D39") "0.02111"
Set-TextValue $ws.Range("E39") "  +1.74%  "
Set-TextValue $ws.Range("E40") "  +3.27%  "
Set-TextValue $ws.Range("D41") "0.5757"
Set-TextValue $ws.Range("E41") "  +2.46%  "
Set-TextValue $ws.Range("D42") "0.1814"
Set-TextValue $ws.Range("E42") "  +1.80%  "
Set-TextValue $ws.Range("D43") "10.02"
Set-TextValue $ws.Range("E43") "  +2.62%  "
Set-TextValue $ws.Range("D44") "12.05"
Set-TextValue $ws.Range("E44") "  +2.41%  "
Set-TextValue $ws.Range("D45") "0.5431"
Set-TextValue $ws.Range("E45") "  +2.64%  "
Set-TextValue $ws.Range("E46") "  +3.30%  "
Set-TextValue $ws.Range("D47") "0.07094"
Set-TextValue $ws.Range("E47") "  -1.06%  "
Set-TextValue $ws.Range("D48") "1.890"
Set-TextValue $ws.Range("E48") "  +3.23%  "
Set-TextValue $ws.Range("D49") "2.577"
Set-TextValue $ws.Range("E49") "  +6.75%  "
Set-TextValue $ws.Range("E50") "  -0.56%  "
Set-TextValue $ws.Range("D51") "1.072"
Set-TextValue $ws.Range("E51") "  -5.39%  "
